# The commit swaps the contents of ppt/theme/theme1.xml (the "Integral" /
# "Red Violet" theme, used by the slide master and therefore by every
# slide) and ppt/theme/theme2.xml (the "Office Theme", used by the notes
# master): after the edit theme1.xml holds the Office Theme palette and
# theme2.xml holds the former Integral/Red Violet palette.
#
# The <a:fmtScheme> (fills/lines/effects) and <a:fontScheme> blocks of the
# two themes are already byte-identical, so the only practical difference
# between "Integral" and "Office Theme" is the 12 colour-scheme entries
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). We reproduce that by
# writing the Office Theme RGB values into the presentation's theme colour
# scheme via the PowerPoint object model.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Index -> (scheme slot, target "Office Theme" RGB as 0xBBGGRR OLE colour)
# 1  dk1       000000 -> 0
# 2  lt1       FFFFFF -> 16777215
# 3  dk2       44546A -> 6968388
# 4  lt2       E7E6E6 -> 15132391
# 5  accent1   5B9BD5 -> 13998939
# 6  accent2   ED7D31 -> 3243501
# 7  accent3   A5A5A5 -> 10855845
# 8  accent4   FFC000 -> 49407
# 9  accent5   4472C4 -> 12874308
# 10 accent6   70AD47 -> 4697456
# 11 hlink     0563C1 -> 12673797
# 12 folHlink  954F72 -> 7491477
$officeThemeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
